# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# cells that get refreshed each time the handback status report is generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 19:07:34"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 19:07:29"
$wsZhCn.Range("K2").Value = "2016-08-25 19:07:46"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 19:07:34"
$wsDeDe.Range("K2").Value = "2016-08-25 19:07:53"
